# Adds Hawkeye simulation results (and completes the OPTGen row) for the
# calculix benchmark on both the Config1 and Config2 sheets.
#
# Context: rows 11-12 on each sheet erroneously held a *duplicate* copy of
# what should have been the calculix/LRU and calculix/SRRIP rows (they were
# mislabeled under cactusADM). That stray data is removed here, and the real
# calculix data (LRU row 15, SRRIP row 16 - which previously only produced
# #DIV/0! placeholders because C:G were empty) is filled in, with brand new
# Hawkeye (row 17) and OPTGen (row 18) results added alongside it.

function Fill-Calculix($ws, $c15, $d15, $e15, $f15, $g15, $c16, $d16, $e16, $f16, $g16, $c17, $d17, $e17, $f17, $g17, $c18, $d18, $e18, $f18) {

    # --- remove the stray/duplicated data that used to sit in rows 11-14 ---
    $null = $ws.Range("C11:I11").ClearContents()
    $null = $ws.Range("C12:I12").ClearContents()
    $null = $ws.Range("H13:I13").ClearContents()
    $null = $ws.Range("H14:I14").ClearContents()

    # --- row 15: calculix / LRU ---
    $ws.Range("C15").Value = $c15
    $ws.Range("D15").Value = $d15
    $ws.Range("E15").Value = $e15
    $ws.Range("F15").Value = $f15
    $ws.Range("G15").Value = $g15
    $ws.Range("H15").Formula = "=(C15/D15)"
    $ws.Range("I15").Formula = "=G15/(C15/1000)"

    # --- row 16: calculix / SRRIP ---
    $ws.Range("C16").Value = $c16
    $ws.Range("D16").Value = $d16
    $ws.Range("E16").Value = $e16
    $ws.Range("F16").Value = $f16
    $ws.Range("G16").Value = $g16
    $ws.Range("H16").Formula = "=(C16/D16)"
    $ws.Range("I16").Formula = "=G16/(C16/1000)"

    # --- row 17: calculix / Hawkeye (new) ---
    $ws.Range("C17").Value = $c17
    $ws.Range("D17").Value = $d17
    $ws.Range("E17").Value = $e17
    $ws.Range("F17").Value = $f17
    $ws.Range("G17").Value = $g17
    $ws.Range("H17").Formula = "=(C17/D17)"
    $ws.Range("I17").Formula = "=G17/(C17/1000)"

    # --- row 18: calculix / OPTGen (new) ---
    $ws.Range("C18").Value = $c18
    $ws.Range("D18").Value = $d18
    $ws.Range("E18").Value = $e18
    $ws.Range("F18").Value = $f18
    $ws.Range("G18").Formula = "=E18-F18"
    $ws.Range("H18").Formula = "=(C18/D18)"
    $ws.Range("I18").Formula = "=G18/(C18/1000)"
    $ws.Range("J18").Formula = "=F18/E18"

    # the active selection moves to C19 after this block of edits
    $null = $ws.Range("C19").Select()
}

$wb = $excel.ActiveWorkbook

$wsConfig1 = $wb.Worksheets.Item("Config1")
Fill-Calculix $wsConfig1 50000001 31836870 405396 341068 64328 50000001 31741436 405396 351651 63745 50000001 32754975 405396 326734 78662 50000001 32754975 12640 10695

$wsConfig2 = $wb.Worksheets.Item("Config2")
Fill-Calculix $wsConfig2 50000000 24198526 425231 359739 65492 50000000 24173346 425232 360877 64355 50000000 24656481 425232 341733 83499 50000000 24656481 2148 7035

$null = $wsConfig1.Activate()
